$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.116.17'
$ws.Range("E2").Value = '  -2.49%  '

$ws.Range("D3").Value = '3.006.09'
$ws.Range("E3").Value = '  -5.26%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.49'
$ws.Range("E5").Value = '  -3.76%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.65'
$ws.Range("E6").Value = '  -5.74%  '

$ws.Range("D8").Value = '3.003.11'
$ws.Range("E8").Value = '  -5.28%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.495'
$ws.Range("E9").Value = '  -2.67%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.134'
$ws.Range("E10").Value = '  -5.90%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.17'
$ws.Range("E11").Value = '  -2.37%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.430'
$ws.Range("E12").Value = '  -5.79%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000223'
$ws.Range("E13").Value = '  -5.53%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.79'
$ws.Range("E14").Value = '  -2.98%  '

$ws.Range("E15").Value = '  -0.58%  '

$ws.Range("D16").Value = '3.506.56'
$ws.Range("E16").Value = '  -5.11%  '

$ws.Range("D17").Value = '61.190.96'
$ws.Range("E17").Value = '  -2.36%  '

$ws.Range("D18").Value = '3.008.35'
$ws.Range("E18").Value = '  -5.14%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.21'
$ws.Range("E19").Value = '  -5.59%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '438.94'
$ws.Range("E20").Value = '  -3.63%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.14'
$ws.Range("E21").Value = '  -6.19%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.661'
$ws.Range("E22").Value = '  -6.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.12'
$ws.Range("E23").Value = '  -6.83%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.54'
$ws.Range("E24").Value = '  -6.61%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '78.82'
$ws.Range("E25").Value = '  -5.97%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.13%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.49'
$ws.Range("E28").Value = '  -7.37%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.22'
$ws.Range("E29").Value = '  -7.11%  '

$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.21'
$ws.Range("E30").Value = '  -8.45%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.87'
$ws.Range("E31").Value = '  -7.60%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '25.45'
$ws.Range("E32").Value = '  -7.05%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0939'
$ws.Range("E33").Value = '  -9.06%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.27'
$ws.Range("E34").Value = '  -5.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.951'
$ws.Range("E35").Value = '  -8.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.56'
$ws.Range("E36").Value = '  -5.24%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '50.03'
$ws.Range("E37").Value = '  -2.33%  '

$ws.Range("D38").Value = '0.0₃0682'
$ws.Range("E38").Value = '  -4.04%  '

$ws.Range("E39").Value = '  -6.16%  '

$ws.Range("B40").Value = 'Cosmos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.73'
$ws.Range("E40").Value = '  -4.06%  '

$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.108'
$ws.Range("E41").Value = '  -3.36%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '370.90'
$ws.Range("E42").Value = '  -8.75%  '

$ws.Range("D43").Value = '2.645.39'
$ws.Range("E43").Value = '  -5.48%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.43'
$ws.Range("E44").Value = '  -10.40%  '

$ws.Range("E45").Value = '  -0.02%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.235'
$ws.Range("E46").Value = '  -6.26%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '119.25'
$ws.Range("E47").Value = '  -5.57%  '

$ws.Range("B48").Value = 'Fetch.AI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.96'
$ws.Range("E48").Value = '  -8.63%  '

$ws.Range("B49").Value = 'Arweave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '32.91'
$ws.Range("E49").Value = '  -5.79%  '

$ws.Range("E50").Value = '  -4.53%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.48'
$ws.Range("E51").Value = '  -8.02%  '
